$wb = $excel.ActiveWorkbook

# Sheet1: update selection from D12 to C12 (this also clears tabSelected
# on Sheet1 once Sheet2 is activated below)
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("C12").Select()

# Sheet2: add a new "Status" column (E) with a "Done" value on row 6,
# make Sheet2 the active tab, and move the selection to E1
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("E1").Value = "Status"
$ws2.Range("E6").Value = "Done"

$ws2.Activate()
$ws2.Range("E1").Select()
